# Updates Fgf2-Gpc4.xlsx LR-pair TPM-derived metrics (columns G:T, rows 2-21)
# to reflect values recomputed with the new TPM input, per commit
# "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.472738
$ws.Cells.Item(2, 8).Value = 1.418214
$ws.Cells.Item(2, 9).Value = 0.0327564895931267
$ws.Cells.Item(2, 10).Value = 0.03397138804734427
$ws.Cells.Item(2, 13).Value = 8.236601666666667
$ws.Cells.Item(2, 14).Value = 24.709805
$ws.Cells.Item(2, 15).Value = 0.1658794346531842
$ws.Cells.Item(2, 16).Value = 0.1779826314087614
$ws.Cells.Item(2, 17).Value = 3.893754598696667
$ws.Cells.Item(2, 18).Value = 35.04379138827
$ws.Cells.Item(2, 19).Value = 0.005433627974930768
$ws.Cells.Item(2, 20).Value = 0.00604631703727448

# Row 3
$ws.Cells.Item(3, 7).Value = 0.472738
$ws.Cells.Item(3, 8).Value = 1.418214
$ws.Cells.Item(3, 9).Value = 0.0327564895931267
$ws.Cells.Item(3, 10).Value = 0.03397138804734427
$ws.Cells.Item(3, 15).Value = 0.6272020870120292
$ws.Cells.Item(3, 16).Value = 0.6729651454676275
$ws.Cells.Item(3, 17).Value = 14.722566517792
$ws.Cells.Item(3, 18).Value = 132.503098660128
$ws.Cells.Item(3, 19).Value = 0.02054493863599688
$ws.Cells.Item(3, 20).Value = 0.02286156009901826

# Row 4
$ws.Cells.Item(4, 7).Value = 0.472738
$ws.Cells.Item(4, 8).Value = 1.418214
$ws.Cells.Item(4, 9).Value = 0.0327564895931267
$ws.Cells.Item(4, 10).Value = 0.03397138804734427
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1445983333333334
$ws.Cells.Item(4, 14).Value = 0.433795
$ws.Cells.Item(4, 15).Value = 0.002912109964258238
$ws.Cells.Item(4, 16).Value = 0.003124588623502439
$ws.Cells.Item(4, 17).Value = 0.06835712690333336
$ws.Cells.Item(4, 18).Value = 0.6152141421300001
$ws.Cells.Item(4, 19).Value = 0.00009539049973826555
$ws.Cells.Item(4, 20).Value = 0.0001061466126173186

# Row 5
$ws.Cells.Item(5, 7).Value = 0.472738
$ws.Cells.Item(5, 8).Value = 1.418214
$ws.Cells.Item(5, 9).Value = 0.0327564895931267
$ws.Cells.Item(5, 10).Value = 0.03397138804734427
$ws.Cells.Item(5, 13).Value = 10.129762
$ws.Cells.Item(5, 14).Value = 20.259524
$ws.Cells.Item(5, 15).Value = 0.2040063683705284
$ws.Cells.Item(5, 16).Value = 0.1459276345001086
$ws.Cells.Item(5, 17).Value = 4.788723428356001
$ws.Cells.Item(5, 18).Value = 28.732340570136
$ws.Cells.Item(5, 19).Value = 0.006682532482460785
$ws.Cells.Item(5, 20).Value = 0.004957364298434213

# Row 6
$ws.Cells.Item(6, 9).Value = 0.822180234441485
$ws.Cells.Item(6, 10).Value = 0.8526739017519405
$ws.Cells.Item(6, 13).Value = 8.236601666666667
$ws.Cells.Item(6, 14).Value = 24.709805
$ws.Cells.Item(6, 15).Value = 0.1658794346531842
$ws.Cells.Item(6, 16).Value = 0.1779826314087614
$ws.Cells.Item(6, 17).Value = 97.73233055735557
$ws.Cells.Item(6, 18).Value = 879.5909750162
$ws.Cells.Item(6, 19).Value = 0.1363827924721759
$ws.Cells.Item(6, 20).Value = 0.1517611447673861

# Row 7
$ws.Cells.Item(7, 9).Value = 0.822180234441485
$ws.Cells.Item(7, 10).Value = 0.8526739017519405
$ws.Cells.Item(7, 15).Value = 0.6272020870120292
$ws.Cells.Item(7, 16).Value = 0.6729651454676275
$ws.Cells.Item(7, 19).Value = 0.5156731589417388
$ws.Cells.Item(7, 20).Value = 0.5738198163289442

# Row 8
$ws.Cells.Item(8, 9).Value = 0.822180234441485
$ws.Cells.Item(8, 10).Value = 0.8526739017519405
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1445983333333334
$ws.Cells.Item(8, 14).Value = 0.433795
$ws.Cells.Item(8, 15).Value = 0.002912109964258238
$ws.Cells.Item(8, 16).Value = 0.003124588623502439
$ws.Cells.Item(8, 17).Value = 1.715747911977778
$ws.Cells.Item(8, 18).Value = 15.4417312078
$ws.Cells.Item(8, 19).Value = 0.002394279253133223
$ws.Cells.Item(8, 20).Value = 0.002664255172971549

# Row 9
$ws.Cells.Item(9, 9).Value = 0.822180234441485
$ws.Cells.Item(9, 10).Value = 0.8526739017519405
$ws.Cells.Item(9, 13).Value = 10.129762
$ws.Cells.Item(9, 14).Value = 20.259524
$ws.Cells.Item(9, 15).Value = 0.2040063683705284
$ws.Cells.Item(9, 16).Value = 0.1459276345001086
$ws.Cells.Item(9, 17).Value = 120.1958390506933
$ws.Cells.Item(9, 18).Value = 721.17503430416
$ws.Cells.Item(9, 19).Value = 0.167730003774437
$ws.Cells.Item(9, 20).Value = 0.1244286854826387

# Row 10
$ws.Cells.Item(10, 7).Value = 0.37892
$ws.Cells.Item(10, 8).Value = 1.13676
$ws.Cells.Item(10, 9).Value = 0.02625574638939025
$ws.Cells.Item(10, 10).Value = 0.02722954016579943
$ws.Cells.Item(10, 13).Value = 8.236601666666667
$ws.Cells.Item(10, 14).Value = 24.709805
$ws.Cells.Item(10, 15).Value = 0.1658794346531842
$ws.Cells.Item(10, 16).Value = 0.1779826314087614
$ws.Cells.Item(10, 17).Value = 3.121013103533333
$ws.Cells.Item(10, 18).Value = 28.0891179318
$ws.Cells.Item(10, 19).Value = 0.004355288367469436
$ws.Cells.Item(10, 20).Value = 0.004846385210759545

# Row 11
$ws.Cells.Item(11, 7).Value = 0.37892
$ws.Cells.Item(11, 8).Value = 1.13676
$ws.Cells.Item(11, 9).Value = 0.02625574638939025
$ws.Cells.Item(11, 10).Value = 0.02722954016579943
$ws.Cells.Item(11, 15).Value = 0.6272020870120292
$ws.Cells.Item(11, 16).Value = 0.6729651454676275
$ws.Cells.Item(11, 17).Value = 11.80077528128
$ws.Cells.Item(11, 18).Value = 106.20697753152
$ws.Cells.Item(11, 19).Value = 0.01646765893148412
$ws.Cells.Item(11, 20).Value = 0.01832453145869382

# Row 12
$ws.Cells.Item(12, 7).Value = 0.37892
$ws.Cells.Item(12, 8).Value = 1.13676
$ws.Cells.Item(12, 9).Value = 0.02625574638939025
$ws.Cells.Item(12, 10).Value = 0.02722954016579943
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1445983333333334
$ws.Cells.Item(12, 14).Value = 0.433795
$ws.Cells.Item(12, 15).Value = 0.002912109964258238
$ws.Cells.Item(12, 16).Value = 0.003124588623502439
$ws.Cells.Item(12, 17).Value = 0.05479120046666668
$ws.Cells.Item(12, 18).Value = 0.4931208042
$ws.Cells.Item(12, 19).Value = 0.0000764596206795806
$ws.Cells.Item(12, 20).Value = 0.00008508111142525961

# Row 13
$ws.Cells.Item(13, 7).Value = 0.37892
$ws.Cells.Item(13, 8).Value = 1.13676
$ws.Cells.Item(13, 9).Value = 0.02625574638939025
$ws.Cells.Item(13, 10).Value = 0.02722954016579943
$ws.Cells.Item(13, 13).Value = 10.129762
$ws.Cells.Item(13, 14).Value = 20.259524
$ws.Cells.Item(13, 15).Value = 0.2040063683705284
$ws.Cells.Item(13, 16).Value = 0.1459276345001086
$ws.Cells.Item(13, 17).Value = 3.83836941704
$ws.Cells.Item(13, 18).Value = 23.03021650224
$ws.Cells.Item(13, 19).Value = 0.005356339469757118
$ws.Cells.Item(13, 20).Value = 0.003973542384920805

# Row 14
$ws.Cells.Item(14, 7).Value = 1.548357
$ws.Cells.Item(14, 8).Value = 3.096714
$ws.Cells.Item(14, 9).Value = 0.1072872076222874
$ws.Cells.Item(14, 10).Value = 0.0741775733180209
$ws.Cells.Item(14, 13).Value = 8.236601666666667
$ws.Cells.Item(14, 14).Value = 24.709805
$ws.Cells.Item(14, 15).Value = 0.1658794346531842
$ws.Cells.Item(14, 16).Value = 0.1779826314087614
$ws.Cells.Item(14, 17).Value = 12.753199846795
$ws.Cells.Item(14, 18).Value = 76.51919908077001
$ws.Cells.Item(14, 19).Value = 0.01779674134590382
$ws.Cells.Item(14, 20).Value = 0.01320231969065769

# Row 15
$ws.Cells.Item(15, 7).Value = 1.548357
$ws.Cells.Item(15, 8).Value = 3.096714
$ws.Cells.Item(15, 9).Value = 0.1072872076222874
$ws.Cells.Item(15, 10).Value = 0.0741775733180209
$ws.Cells.Item(15, 15).Value = 0.6272020870120292
$ws.Cells.Item(15, 16).Value = 0.6729651454676275
$ws.Cells.Item(15, 17).Value = 48.22076694868801
$ws.Cells.Item(15, 18).Value = 289.3246016921281
$ws.Cells.Item(15, 19).Value = 0.06729076053039151
$ws.Cells.Item(15, 20).Value = 0.04991892141839754

# Row 16
$ws.Cells.Item(16, 7).Value = 1.548357
$ws.Cells.Item(16, 8).Value = 3.096714
$ws.Cells.Item(16, 9).Value = 0.1072872076222874
$ws.Cells.Item(16, 10).Value = 0.0741775733180209
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1445983333333334
$ws.Cells.Item(16, 14).Value = 0.433795
$ws.Cells.Item(16, 15).Value = 0.002912109964258238
$ws.Cells.Item(16, 16).Value = 0.003124588623502439
$ws.Cells.Item(16, 17).Value = 0.2238898416050001
$ws.Cells.Item(16, 18).Value = 1.34333904963
$ws.Cells.Item(16, 19).Value = 0.0003124321463543054
$ws.Cells.Item(16, 20).Value = 0.0002317744017085061

# Row 17
$ws.Cells.Item(17, 7).Value = 1.548357
$ws.Cells.Item(17, 8).Value = 3.096714
$ws.Cells.Item(17, 9).Value = 0.1072872076222874
$ws.Cells.Item(17, 10).Value = 0.0741775733180209
$ws.Cells.Item(17, 13).Value = 10.129762
$ws.Cells.Item(17, 14).Value = 20.259524
$ws.Cells.Item(17, 15).Value = 0.2040063683705284
$ws.Cells.Item(17, 16).Value = 0.1459276345001086
$ws.Cells.Item(17, 17).Value = 15.684487901034
$ws.Cells.Item(17, 18).Value = 62.73795160413601
$ws.Cells.Item(17, 19).Value = 0.02188727359963771
$ws.Cells.Item(17, 20).Value = 0.01082455780725716

# Row 18
$ws.Cells.Item(18, 7).Value = 0.16626
$ws.Cells.Item(18, 8).Value = 0.49878
$ws.Cells.Item(18, 9).Value = 0.01152032195371061
$ws.Cells.Item(18, 10).Value = 0.01194759671689489
$ws.Cells.Item(18, 13).Value = 8.236601666666667
$ws.Cells.Item(18, 14).Value = 24.709805
$ws.Cells.Item(18, 15).Value = 0.1658794346531842
$ws.Cells.Item(18, 16).Value = 0.1779826314087614
$ws.Cells.Item(18, 17).Value = 1.3694173931
$ws.Cells.Item(18, 18).Value = 12.3247565379
$ws.Cells.Item(18, 19).Value = 0.001910984492704182
$ws.Cells.Item(18, 20).Value = 0.002126464702683632

# Row 19
$ws.Cells.Item(19, 7).Value = 0.16626
$ws.Cells.Item(19, 8).Value = 0.49878
$ws.Cells.Item(19, 9).Value = 0.01152032195371061
$ws.Cells.Item(19, 10).Value = 0.01194759671689489
$ws.Cells.Item(19, 15).Value = 0.6272020870120292
$ws.Cells.Item(19, 16).Value = 0.6729651454676275
$ws.Cells.Item(19, 17).Value = 5.17786577184
$ws.Cells.Item(19, 18).Value = 46.60079194656
$ws.Cells.Item(19, 19).Value = 0.00722556997241779
$ws.Cells.Item(19, 20).Value = 0.00804031616257372

# Row 20
$ws.Cells.Item(20, 7).Value = 0.16626
$ws.Cells.Item(20, 8).Value = 0.49878
$ws.Cells.Item(20, 9).Value = 0.01152032195371061
$ws.Cells.Item(20, 10).Value = 0.01194759671689489
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 0.6666666666666666
$ws.Cells.Item(20, 13).Value = 0.1445983333333334
$ws.Cells.Item(20, 14).Value = 0.433795
$ws.Cells.Item(20, 15).Value = 0.002912109964258238
$ws.Cells.Item(20, 16).Value = 0.003124588623502439
$ws.Cells.Item(20, 17).Value = 0.0240409189
$ws.Cells.Item(20, 18).Value = 0.2163682701
$ws.Cells.Item(20, 19).Value = 0.00003354844435286359
$ws.Cells.Item(20, 20).Value = 0.00003733132477980487

# Row 21
$ws.Cells.Item(21, 7).Value = 0.16626
$ws.Cells.Item(21, 8).Value = 0.49878
$ws.Cells.Item(21, 9).Value = 0.01152032195371061
$ws.Cells.Item(21, 10).Value = 0.01194759671689489
$ws.Cells.Item(21, 13).Value = 10.129762
$ws.Cells.Item(21, 14).Value = 20.259524
$ws.Cells.Item(21, 15).Value = 0.2040063683705284
$ws.Cells.Item(21, 16).Value = 0.1459276345001086
$ws.Cells.Item(21, 17).Value = 1.68417423012
$ws.Cells.Item(21, 18).Value = 10.10504538072
$ws.Cells.Item(21, 19).Value = 0.002350219044235771
$ws.Cells.Item(21, 20).Value = 0.001743484526857735

Write-Output "Updated Fgf2-Gpc4 TPM values (230 cells)"
